$p = $ppt.ActivePresentation

$newStyle = "{3693021D-87AE-4348-86FD-531A93FFC3C5}"

foreach ($idx in 14, 15, 16) {
    $s = $p.Slides.Item($idx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newStyle)
        }
    }
}
